$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# E1 was CA (EXPIRATIONDATE column header actually held "EFFECTIVEDATE" text before;
# after the edit E1 should read "EFFECTIVEDATE" still, but now referencing the new
# shared-string ordering). G1 should now read "MSRP_VERSION" (unchanged text) but via a
# different shared-string slot; we just need the cell to end up with the right text.
$ws.Range("E1").Value = "EFFECTIVEDATE"
$ws.Range("G1").Value = "MSRP_VERSION"

# --- Update row 2 (existing data row) ---
$ws.Range("A2").Value = "AAA_CSA"
$ws.Range("B2").Value = "SELECT"
$ws.Range("C2").Value = "CA"
$ws.Range("D2").Value = "SYMBOL_2000"
$ws.Range("E2").Value = 20000101
$ws.Range("F2").Value = 20190701

# --- Add new row 3 ---
$ws.Range("A3").Value = "AAA_CSA"
$ws.Range("B3").Value = "SELECT"
$ws.Range("C3").Value = "CA"
# Introduce the new shared strings in the same order the original workbook used
# (SYMBOL_2018 first, then MSRP_2000_SELECT) so the rebuilt sharedStrings table
# lines up with the target ordering.
$ws.Range("D3").Value = "SYMBOL_2018"
$ws.Range("G2").Value = "MSRP_2000_SELECT"
$ws.Range("E3").Value = 20190702
$ws.Range("F3").Value = 99999999
$ws.Range("G3").Value = "MSRP_2000_SELECT"

# Apply the new "bordered / highlighted" style to F2 (new custom format: Good font,
# light fill, thin black border all around), then copy that exact formatting onto E3
# so both cells end up sharing a single new cell style record.
$f2 = $ws.Range("F2")
$f2.Style = "Good"
$f2.Interior.ThemeColor = 2
$f2.Interior.TintAndShade = 0
$f2.Borders.LineStyle = 1
$f2.Borders.Weight = 2

$f2.Copy()
$ws.Range("E3").PasteSpecial(-4122)

# --- Sheet view / dimension bookkeeping ---
$ws.Range("G9").Select()
